# Fruta / hortaliza, semanal
# Insert a new weekly price record as row 141 (pushing the former rows 141-167
# down to 142-168), for "Feria Lagunitas de Puerto Montt - Haba".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 141..167 down to 142..168 by inserting a new row at 141.
$ws.Rows.Item(141).Insert()

# Populate the newly inserted row 141 with the new weekly record.
$ws.Range("A141").Value = 4
$ws.Range("B141").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C141").Value = "Los Lagos"
$ws.Range("D141").Value = 45209
$ws.Range("E141").Value = 10
$ws.Range("F141").Value = 100112026
$ws.Range("G141").Value = "Haba"
$ws.Range("H141").Value = "Sin especificar"
$ws.Range("I141").Value = "Primera"
$ws.Range("J141").Value = 120
$ws.Range("K141").Value = 18000
$ws.Range("L141").Value = 18000
$ws.Range("M141").Value = 18000
$ws.Range("N141").Value = "$/saco 25 kilos"
$ws.Range("O141").Value = "Provincia de Limarí"
$ws.Range("P141").Value = 720
$ws.Range("Q141").Value = 25
$ws.Range("R141").Value = "Hortaliza"
